# Insert a new weekly price record at row 108 for
# "Terminal La Palmera de La Serena" / Jengibre, pushing the existing
# rows 108:175 down to 109:176 (dimension grows from A1:R175 to A1:R176).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 108 - this shifts every
# row from 108 down to 175 one row lower (now 109..176) and inherits
# the row's existing formatting (date style on column D).
$ws.Rows.Item(108).Insert()

# Populate the newly inserted row 108 with the new record.
$ws.Cells.Item(108, 1).Value  = 8
$ws.Cells.Item(108, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(108, 3).Value  = "Coquimbo"
$ws.Cells.Item(108, 4).Value  = 45161
$ws.Cells.Item(108, 5).Value  = 4
$ws.Cells.Item(108, 6).Value  = 100114007
$ws.Cells.Item(108, 7).Value  = "Jengibre"
$ws.Cells.Item(108, 8).Value  = "Sin especificar"
$ws.Cells.Item(108, 9).Value  = "Primera"
$ws.Cells.Item(108, 10).Value = 420
$ws.Cells.Item(108, 11).Value = 18000
$ws.Cells.Item(108, 12).Value = 19000
$ws.Cells.Item(108, 13).Value = 18500
$ws.Cells.Item(108, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(108, 15).Value = "Perú"
$ws.Cells.Item(108, 16).Value = 1423
$ws.Cells.Item(108, 17).Value = 13
$ws.Cells.Item(108, 18).Value = "Hortaliza"
